# The deck ships with two theme parts:
#   ppt/theme/theme1.xml -> linked from the notes master ("Office Theme")
#   ppt/theme/theme2.xml -> linked from the slide master / main deck ("Integral" / "Red Violet")
#
# The authored edit swaps the two themes' contents so the slide master
# (the theme that actually paints the slides) now uses the plain "Office"
# color palette instead of the pink/violet "Integral" palette. Font scheme
# and format (fill/line/effect) scheme are identical between the two themes,
# so the only real delta is the 12 theme colors (dk1/lt1/dk2/lt2/accent1-6/
# hlink/folHlink). Re-point the live theme's color scheme to the "Office"
# values via the master's ThemeColorScheme, color by color.

$p  = $ppt.ActivePresentation
$sm = $p.SlideMaster
$tcs = $sm.Theme.ThemeColorScheme

# index : role      : target "Office" RGB (stored as BGR-packed long for COM)
$tcs.Item(1).RGB  = 0         # dk1      000000
$tcs.Item(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388   # dk2      44546A
$tcs.Item(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407     # accent4  FFC000
$tcs.Item(9).RGB  = 12874308  # accent5  4472C4
$tcs.Item(10).RGB = 4697456   # accent6  70AD47
$tcs.Item(11).RGB = 12673797  # hlink    0563C1
$tcs.Item(12).RGB = 7491477   # folHlink 954F72
